$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QUITAFACIL")

# Update TELEFONES for sergio pj (row 5): prepend new phone number
$ws.Range("B5").Value = "55415865, 4254854445, 11982662870"

# NOME stays the same value, just ensure it's intact
$ws.Range("E5").Value = "sergio pj"

# Move the active selection / view to B5 (was E5)
$ws.Range("B5").Select()
